$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.149.03"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "3.531.35"
$ws.Range("E3").Value = "  +2.79%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").Value = "3.528.73"
$ws.Range("E7").Value = "  +2.75%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  -2.21%  "

$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.32%  "

$ws.Range("E12").Value = "  +2.71%  "

$ws.Range("D13").Value = "4.128.18"
$ws.Range("E13").Value = "  +2.82%  "

$ws.Range("E14").Value = "  +2.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.68%  "

$ws.Range("D16").Value = "3.532.19"
$ws.Range("E16").Value = "  +2.79%  "

$ws.Range("E17").Value = "  +1.45%  "

$ws.Range("D18").Value = "65.212.63"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.51%  "

$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("E23").Value = "  +3.38%  "

$ws.Range("D24").Value = "3.670.73"
$ws.Range("E24").Value = "  +2.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.73"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +7.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.56%  "

$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("D32").Value = "3.542.38"
$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.26%  "

$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("E36").Value = "  +7.58%  "

$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.41%  "

$ws.Range("E39").Value = "  +4.76%  "

$ws.Range("E40").Value = "  +3.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0797"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.71%  "

$ws.Range("E42").Value = "  -0.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.34%  "

$ws.Range("E48").Value = "  +5.64%  "

$ws.Range("E49").Value = "  +3.44%  "

$ws.Range("D50").Value = "2.384.73"
$ws.Range("E50").Value = "  +7.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "301.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.63%  "

